$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Inhbb"
$ws.Range("C2").Value = "Acvr2b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 3.345805333333333
$ws.Range("H2").Value = 10.037416
$ws.Range("I2").Value = 0.3489465220682754
$ws.Range("J2").Value = 0.3489465220682754
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5292956666666667
$ws.Range("N2").Value = 1.587887
$ws.Range("O2").Value = 0.135651968140022
$ws.Range("P2").Value = 0.1356519681400219
$ws.Range("Q2").Value = 1.770920264443556
$ws.Range("R2").Value = 15.938282379992
$ws.Range("S2").Value = 0.04733528249417716
$ws.Range("T2").Value = 0.04733528249417715

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Inhbb"
$ws.Range("C3").Value = "Acvr2b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 3.345805333333333
$ws.Range("H3").Value = 10.037416
$ws.Range("I3").Value = 0.3489465220682754
$ws.Range("J3").Value = 0.3489465220682754
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.362890666666667
$ws.Range("N3").Value = 4.088672000000001
$ws.Range("O3").Value = 0.3492921120199358
$ws.Range("P3").Value = 0.3492921120199358
$ws.Range("Q3").Value = 4.559966861283556
$ws.Range("R3").Value = 41.039701751552
$ws.Range("S3").Value = 0.121884267675239
$ws.Range("T3").Value = 0.121884267675239

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Inhbb"
$ws.Range("C4").Value = "Acvr2b"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 3.345805333333333
$ws.Range("H4").Value = 10.037416
$ws.Range("I4").Value = 0.3489465220682754
$ws.Range("J4").Value = 0.3489465220682754
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.009678666666666
$ws.Range("N4").Value = 6.029036
$ws.Range("O4").Value = 0.5150559198400423
$ws.Range("P4").Value = 0.5150559198400423
$ws.Range("Q4").Value = 6.723993601219553
$ws.Range("R4").Value = 60.51594241097599
$ws.Range("S4").Value = 0.1797269718988592
$ws.Range("T4").Value = 0.1797269718988592

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Inhbb"
$ws.Range("C5").Value = "Acvr2b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.480061666666667
$ws.Range("H5").Value = 16.440185
$ws.Range("I5").Value = 0.5715360783999618
$ws.Range("J5").Value = 0.5715360783999618
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5292956666666667
$ws.Range("N5").Value = 1.587887
$ws.Range("O5").Value = 0.135651968140022
$ws.Range("P5").Value = 0.1356519681400219
$ws.Range("Q5").Value = 2.900572893232778
$ws.Range("R5").Value = 26.105156039095
$ws.Range("S5").Value = 0.07752999389798471
$ws.Range("T5").Value = 0.07752999389798469

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Inhbb"
$ws.Range("C6").Value = "Acvr2b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.480061666666667
$ws.Range("H6").Value = 16.440185
$ws.Range("I6").Value = 0.5715360783999618
$ws.Range("J6").Value = 0.5715360783999618
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.362890666666667
$ws.Range("N6").Value = 4.088672000000001
$ws.Range("O6").Value = 0.3492921120199358
$ws.Range("P6").Value = 0.3492921120199358
$ws.Range("Q6").Value = 7.468724898257779
$ws.Range("R6").Value = 67.21852408432001
$ws.Range("S6").Value = 0.1996330439199143
$ws.Range("T6").Value = 0.1996330439199143

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Inhbb"
$ws.Range("C7").Value = "Acvr2b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.480061666666667
$ws.Range("H7").Value = 16.440185
$ws.Range("I7").Value = 0.5715360783999618
$ws.Range("J7").Value = 0.5715360783999618
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.009678666666666
$ws.Range("N7").Value = 6.029036
$ws.Range("O7").Value = 0.5150559198400423
$ws.Range("P7").Value = 0.5150559198400423
$ws.Range("Q7").Value = 11.01316302351778
$ws.Range("R7").Value = 99.11846721165999
$ws.Range("S7").Value = 0.2943730405820628
$ws.Range("T7").Value = 0.2943730405820628

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Inhbb"
$ws.Range("C8").Value = "Acvr2b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7624369999999999
$ws.Range("H8").Value = 2.287311
$ws.Range("I8").Value = 0.07951739953176286
$ws.Range("J8").Value = 0.07951739953176286
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5292956666666667
$ws.Range("N8").Value = 1.587887
$ws.Range("O8").Value = 0.135651968140022
$ws.Range("P8").Value = 0.1356519681400219
$ws.Range("Q8").Value = 0.4035546002063333
$ws.Range("R8").Value = 3.631991401857
$ws.Range("S8").Value = 0.01078669174786009
$ws.Range("T8").Value = 0.01078669174786009

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Inhbb"
$ws.Range("C9").Value = "Acvr2b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7624369999999999
$ws.Range("H9").Value = 2.287311
$ws.Range("I9").Value = 0.07951739953176286
$ws.Range("J9").Value = 0.07951739953176286
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.362890666666667
$ws.Range("N9").Value = 4.088672000000001
$ws.Range("O9").Value = 0.3492921120199358
$ws.Range("P9").Value = 0.3492921120199358
$ws.Range("Q9").Value = 1.039118271221333
$ws.Range("R9").Value = 9.352064440992001
$ws.Range("S9").Value = 0.0277748004247825
$ws.Range("T9").Value = 0.0277748004247825

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Inhbb"
$ws.Range("C10").Value = "Acvr2b"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7624369999999999
$ws.Range("H10").Value = 2.287311
$ws.Range("I10").Value = 0.07951739953176286
$ws.Range("J10").Value = 0.07951739953176286
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.009678666666666
$ws.Range("N10").Value = 6.029036
$ws.Range("O10").Value = 0.5150559198400423
$ws.Range("P10").Value = 0.5150559198400423
$ws.Range("Q10").Value = 1.532253373577333
$ws.Range("R10").Value = 13.790280362196
$ws.Range("S10").Value = 0.04095590735912027
$ws.Range("T10").Value = 0.04095590735912027
